# Auto-generated Excel COM-interop edit script
$wb = $excel.ActiveWorkbook

# --- Resolve the two worksheets by their CURRENT (pre-edit) names ---
$wsApk = $wb.Worksheets.Item("TitleStringTable")   # becomes "InApkStringTable"
$wsStr = $wb.Worksheets.Item("StringTable")        # stays "StringTable"

# --- Rename the first sheet ---
$wsApk.Name = "InApkStringTable"

# --- Clear existing contents/formatting of both sheets before rebuilding them ---
$wsApk.Cells.Clear()
$wsStr.Cells.Clear()

# --- InApkStringTable (sheet1, becomes the big / tab-selected sheet) ---
$wsApk.Cells.Item(1,1).Value = "id|String"
$wsApk.Cells.Item(1,2).Value = "kor|String"
$wsApk.Cells.Item(1,3).Value = "eng|String"

$wsApk.Cells.Item(2,1).Value = "CommonUI_Ok"
$wsApk.Cells.Item(2,2).Value = "확인"
$wsApk.Cells.Item(2,3).Value = "OK"

$wsApk.Cells.Item(3,1).Value = "CommonUI_Cancel"
$wsApk.Cells.Item(3,2).Value = "취소"
$wsApk.Cells.Item(3,3).Value = "Cancel"

$wsApk.Cells.Item(4,1).Value = "CommonUI_Yes"
$wsApk.Cells.Item(4,2).Value = "예"
$wsApk.Cells.Item(4,3).Value = "Yes"

$wsApk.Cells.Item(5,1).Value = "CommonUI_No"
$wsApk.Cells.Item(5,2).Value = "아니오"
$wsApk.Cells.Item(5,3).Value = "No"

$wsApk.Cells.Item(6,1).Value = "SystemUI_NeedPatch"
$wsApk.Cells.Item(6,2).Value = "새로운 패치가 있습니다.`n<color=cyan>와이파이</color>를 사용해서 다운로드 받으시길 권장합니다."
$wsApk.Cells.Item(6,2).WrapText = $true
$wsApk.Cells.Item(6,3).Value = "New Patch!`nRecommended to download with <color=cyan>WIFI</color>"
$wsApk.Cells.Item(6,3).WrapText = $true

$wsApk.Cells.Item(7,1).Value = "SystemUI_DisconnectServer"
$wsApk.Cells.Item(7,2).Value = "서버와의 접속이 원활하지 않습니다.`n프로그램을 재시작합니다."
$wsApk.Cells.Item(7,2).WrapText = $true
$wsApk.Cells.Item(7,3).Value = "Bad connection`nRestarting the app"
$wsApk.Cells.Item(7,3).WrapText = $true

$wsApk.Cells.Item(8,1).Value = "SystemUI_Mainternance"
$wsApk.Cells.Item(8,2).Value = "서버 점검 중입니다."
$wsApk.Cells.Item(8,2).WrapText = $true
$wsApk.Cells.Item(8,3).Value = "We're on a mainternance"
$wsApk.Cells.Item(8,3).WrapText = $true

$wsApk.Cells.Item(9,1).Value = "SystemUI_MainternanceDetail"
$wsApk.Cells.Item(9,2).Value = "5월 29일 화요일 오전 3시부터 대략 11시까지 정기 점검 예정입니다. 이 동안 서버 및 웹 서비스 등을 사용하실 수 없습니다.`n감사합니다."
$wsApk.Cells.Item(9,2).WrapText = $true
$wsApk.Cells.Item(9,3).Value = "We will be performing scheduled maintenance on Tuesday, May 29th. Maintenance will begin at 3:00 AM PDT and conclude at approximately 11:00 AM PDT. During this time, servers and many web services will be unavailable.`nThank you for your patience."
$wsApk.Cells.Item(9,3).WrapText = $true

$wsApk.Cells.Item(10,1).Value = "SystemUI_Downloading"
$wsApk.Cells.Item(10,2).Value = "다운로드 중"
$wsApk.Cells.Item(10,3).Value = "Downloading"

$wsApk.Cells.Item(11,1).Value = "GameUI_Swappable"
$wsApk.Cells.Item(11,2).Value = "교체 가능"
$wsApk.Cells.Item(11,3).Value = "Can be swapped"

$wsApk.Cells.Item(12,1).Value = "GameUI_TouchToMove"
$wsApk.Cells.Item(12,2).Value = "터치하여 이동하세요"
$wsApk.Cells.Item(12,3).Value = "Touch to move"

$wsApk.Cells.Item(13,1).Value = "GameUI_RepiarPack"
$wsApk.Cells.Item(13,2).Value = "{0}레벨 수리킷 사용"
$wsApk.Cells.Item(13,3).Value = "Use Repair Kit Lv. {0}"

$wsApk.Cells.Item(14,1).Value = "GameUI_Play"
$wsApk.Cells.Item(14,2).Value = "진행"
$wsApk.Cells.Item(14,3).Value = "Play"

$wsApk.Cells.Item(15,1).Value = "GameUI_PossibleAfterTraining"
$wsApk.Cells.Item(15,2).Value = "훈련 챕터 클리어 후 진행 가능"
$wsApk.Cells.Item(15,3).Value = "Possible to play after the training chapter"

$wsApk.Cells.Item(16,1).Value = "GameUI_Shop"
$wsApk.Cells.Item(16,2).Value = "상점"
$wsApk.Cells.Item(16,3).Value = "Shop"

$wsApk.Cells.Item(17,1).Value = "GameUI_UnderConstruction"
$wsApk.Cells.Item(17,2).Value = "개발 중"
$wsApk.Cells.Item(17,3).Value = "Under Construction"

$wsApk.Columns.Item(1).ColumnWidth = 28.857142857142858
$wsApk.Columns.Item(2).ColumnWidth = 96.57142857142857
$wsApk.Columns.Item(3).ColumnWidth = 77.71428571428571

# --- StringTable (sheet2, shrinks down to just the Skill_* rows) ---
$wsStr.Cells.Item(1,1).Value = "id|String"
$wsStr.Cells.Item(1,2).Value = "kor|String"
$wsStr.Cells.Item(1,3).Value = "eng|String"

$wsStr.Cells.Item(2,1).Value = "Skill_ActiveOne001_Name"
$wsStr.Cells.Item(2,2).Value = "하트가 폭발한다"
$wsStr.Cells.Item(2,3).Value = "Hearts exploding!"

$wsStr.Cells.Item(3,1).Value = "Skill_ActiveOne001_Description"
$wsStr.Cells.Item(3,2).Value = "대미지를 <color=#FFFF00>{0}%,{1}%,{2}%,{3}%</color> 먹이고 적이 죽이면 하트를 반드시 떨어뜨린다`n적이 죽지 않으면 기절을 먹인다"
$wsStr.Cells.Item(3,2).WrapText = $true
$wsStr.Cells.Item(3,3).Value = "Deal <color=#FFFF00>{0}%,{1}%,{2}%,{3}%</color> and then the enemy drops a heart if killed. If not, stunned."

$wsStr.Cells.Item(4,1).Value = "Skill_ActiveOne002_Description"
$wsStr.Cells.Item(4,2).Value = "대미지를 {0}% 먹이고 적이 죽이면 하트를 반드시 떨어뜨린다`n적이 죽지 않으면 기절을 먹인다"
$wsStr.Cells.Item(4,3).Value = "Deal <color=#FFFF00>{0}%</color> and then the enemy drops a heart if killed. If not, stunned."

$wsStr.Columns.Item(1).ColumnWidth = 28.857142857142858
$wsStr.Columns.Item(2).ColumnWidth = 96.57142857142857
$wsStr.Columns.Item(3).ColumnWidth = 77.71428571428571

# --- Row height for the header row (explicit custom height, unaffected by wrap) ---
$wsApk.Rows.Item(1).RowHeight = 27
$wsStr.Rows.Item(1).RowHeight = 27

# --- Make InApkStringTable the active / selected tab ---
$wsApk.Activate()

